{"js": "// The edit replaces the sentence about signing the paper consent form with\n// a sentence about verbal consent (see commit message \"update consent form\n// to verbal\") inside the \"About this consent form\" paragraph.\n\nconst oldText =\n  \"If you decide to participate in this research you will be asked to sign this form. A copy of the signed form will be provided to you for your record.\";\nconst newText =\n  \"If you decide to participate in this research, we will ask for verbal consent. You can save this for consent form for your record.\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + oldText);\n}\n\n// Replace just the trailing part that changed, keeping the unaffected\n// leading text (\"Please read this form carefully. ... at any time. \")\n// untouched by only acting on the matched range.\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The edit replaces the sentence about signing the paper consent form with\n# a sentence about verbal consent (see commit message \"update consent form\n# to verbal\") inside the \"About this consent form\" paragraph.\n\n$d = $word.ActiveDocument\n\n$oldText = \"If you decide to participate in this research you will be asked to sign this form. A copy of the signed form will be provided to you for your record.\"\n$newText = \"If you decide to participate in this research, we will ask for verbal consent. You can save this for consent form for your record.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$find.Execute([ref]$null, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$null, 2) | Out-Null\n"}
